# Issue #20: Printing the output of grep
#
# Adds two new columns ("pattern" renamed from the old "regex" header text,
# plus a brand-new "regex" / "case-sensitive" pair of boolean flag columns),
# and adds several new database/language heuristic rows (HyperSQL, Derby,
# H2, Oracle, PostgreSQL), while also tightening the existing MySQL jdbc
# pattern from "jdbc:mysql:" to "jdbc:mysql".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header row: C1 becomes "pattern", D1/E1 are new "regex" /
#    "case-sensitive" headers. Copy C1's bold+border formatting onto the
#    two new header cells first, then overwrite the text.
# ---------------------------------------------------------------------
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1:E1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Cells.Item(1, 3).Value = "pattern"
$ws.Cells.Item(1, 4).Value = "regex"
$ws.Cells.Item(1, 5).Value = "case-sensitive"

# ---------------------------------------------------------------------
# 2. Data rows: fix the MySQL jdbc pattern ("jdbc:mysql:" -> "jdbc:mysql"),
#    append the new flag columns (D = regex => False, E = case-sensitive
#    => True) to every existing row, and add the new heuristic rows for
#    HyperSQL, Derby, H2, Oracle and PostgreSQL.
# ---------------------------------------------------------------------

# Full target table (rows 2..9), columns A..E.
$rows = @(
    @("Java",  "MySQL",      "jdbc:mysql",                               $false, $true),
    @("Java",  "MySQL",      "mysqlx:",                                  $false, $true),
    @("C#",    "MySQL",      "MySql.Data.MySqlClient.MySqlConnection",   $false, $true),
    @("Java ", "HyperSQL",   "jdbc:hsqldb",                              $false, $true),
    @("Java",  "Derby",      "jdbc:derby",                               $false, $true),
    @("Java",  "H2",         "jdbc:h2",                                  $false, $true),
    @("Java",  "Oracle",     "jdbc:oracle",                              $false, $true),
    @("Java",  "PostgreSQL", "jdbc:postgresql",                          $false, $true)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# ---------------------------------------------------------------------
# 3. Column widths: best-fit all five columns to their new contents.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(4).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(5).EntireColumn.AutoFit() | Out-Null

# ---------------------------------------------------------------------
# 4. Selection moves to C18 in the final state.
# ---------------------------------------------------------------------
$ws.Range("C18").Select() | Out-Null
